# Apply the edits described by the diff:
#  1. Remove the stray empty cell B11 on the "ODI Batting" sheet.
#  2. Add a new worksheet "ODI Batting Extra" (after "ODI Bowling") with
#     MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#     PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.

$wb = $excel.ActiveWorkbook

# --- 1. ODI Batting: drop the leftover empty inline-string cell at B11 ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B11").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet after the last sheet ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Copy the header formatting (bold, centered, bordered) used by the other
# sheets so the new header row reuses the same cell style as the rest of
# the workbook. "ODI Bowling" has 7 header columns, so A1:F1 is fully
# populated there (unlike "Player Info", which only has 4 columns).
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122) # xlPasteFormats

# --- Header row ----------------------------------------------------------
$extra.Cells.Item(1,1).Value = "MATCH_CODE"
$extra.Cells.Item(1,2).Value = "BATTING_POSITION"
$extra.Cells.Item(1,3).Value = "NUM_4"
$extra.Cells.Item(1,4).Value = "NUM_6"
$extra.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Cells.Item(1,6).Value = "MAN_OF_MATCH"

# --- Data rows -------------------------------------------------------------
# Column A (MATCH_CODE) and C/D/E (NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL)
# are text in the source data even though they look numeric, so they are
# entered with a leading apostrophe to force text storage; column B
# (BATTING_POSITION) is a genuine number; column F (MAN_OF_MATCH) is plain
# text.

$extra.Cells.Item(2,1).Value = "'3859"
$extra.Cells.Item(2,2).Value = 4
$extra.Cells.Item(2,3).Value = "'1"
$extra.Cells.Item(2,4).Value = "'0"
$extra.Cells.Item(2,5).Value = "'2.66%"
$extra.Cells.Item(2,6).Value = "NO"

$extra.Cells.Item(3,1).Value = "'3861"
$extra.Cells.Item(3,6).Value = "NO"

$extra.Cells.Item(4,1).Value = "'4375"
$extra.Cells.Item(4,2).Value = 6
$extra.Cells.Item(4,3).Value = "'2"
$extra.Cells.Item(4,4).Value = "'2"
$extra.Cells.Item(4,5).Value = "'10.49%"
$extra.Cells.Item(4,6).Value = "NO"

$extra.Cells.Item(5,1).Value = "'4376"
$extra.Cells.Item(5,2).Value = 6
$extra.Cells.Item(5,3).Value = "'3"
$extra.Cells.Item(5,4).Value = "'1"
$extra.Cells.Item(5,5).Value = "'9.36%"
$extra.Cells.Item(5,6).Value = "NO"

$extra.Cells.Item(6,1).Value = "'4432"
$extra.Cells.Item(6,6).Value = "NO"

$extra.Cells.Item(7,1).Value = "'4433"
$extra.Cells.Item(7,6).Value = "NO"

$extra.Cells.Item(8,1).Value = "'4434"
$extra.Cells.Item(8,2).Value = 6
$extra.Cells.Item(8,3).Value = "'2"
$extra.Cells.Item(8,4).Value = "'0"
$extra.Cells.Item(8,5).Value = "'6.47%"
$extra.Cells.Item(8,6).Value = "NO"

$extra.Cells.Item(9,1).Value = "'4564"
$extra.Cells.Item(9,2).Value = 6
$extra.Cells.Item(9,3).Value = "'0"
$extra.Cells.Item(9,4).Value = "'0"
$extra.Cells.Item(9,5).Value = "'0.89%"
$extra.Cells.Item(9,6).Value = "NO"

$extra.Cells.Item(10,1).Value = "'4565"
$extra.Cells.Item(10,6).Value = "NO"

$extra.Cells.Item(11,1).Value = "'4567"
$extra.Cells.Item(11,2).Value = 6
$extra.Cells.Item(11,6).Value = "NO"

# The leading-apostrophe entries above leave a "quote prefix" cell style on
# A2:A11 / C2:E11; strip that back off so the data rows end up on the
# default (unstyled) cell format, matching the rest of the workbook.
$extra.Range("A2:A11").Style = "Normal"
$extra.Range("C2:E11").Style = "Normal"
